# Generate Report for Handback
# Update the timestamp values recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" for the 8974694f... file
$wsOverview.Range("G2").Value = "2016-08-29 07:06:10"

# zh-cn!H2 - "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-08-29 07:05:59"

# zh-cn!K2 - "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-08-29 07:06:30"

# de-de!H2 - "Correspond Handoff Datetime"
$wsDeDe.Range("H2").Value = "2016-08-29 07:06:10"

# de-de!K2 - "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-08-29 07:06:37"
